# Apply updated crypto price/volume figures (Tue Sep 26 13:46:44 UTC 2023 refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.248.15"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "1.588.80"
$ws.Range("E3").Value = "  +0.91%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'212.18"
$ws.Range("E6").Value = "  +0.84%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  +0.72%  "
$ws.Range("D9").Value = "'0.0608"
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("E10").Value = "  -1.09%  "
$ws.Range("D11").Value = "'0.0849"
$ws.Range("E11").Value = "  +0.69%  "
$ws.Range("D12").Value = "1.811.20"
$ws.Range("E12").Value = "  +0.85%  "
$ws.Range("D13").Value = "1.599.49"
$ws.Range("E13").Value = "  +1.72%  "
$ws.Range("E14").Value = "  -1.16%  "
$ws.Range("D15").Value = "'0.520"
$ws.Range("E15").Value = "  +1.13%  "
$ws.Range("D16").Value = "'64.26"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").Value = "26.238.47"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").Value = "0.0₃0727"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("E19").Value = "  +1.54%  "
$ws.Range("D20").Value = "'213.05"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").Value = "'4.27"
$ws.Range("E22").Value = "  +0.29%  "
$ws.Range("D23").Value = "'2.17"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").Value = "'9.00"
$ws.Range("E24").Value = "  +2.01%  "
$ws.Range("D25").Value = "'143.54"
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").Value = "'7.01"
$ws.Range("E27").Value = "  +0.48%  "
$ws.Range("E28").Value = "  -0.66%  "
$ws.Range("D29").Value = "'15.18"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("E30").Value = "  -1.88%  "
$ws.Range("E31").Value = "  +1.45%  "
$ws.Range("E32").Value = "  -0.23%  "
$ws.Range("D33").Value = "1.338.66"
$ws.Range("E33").Value = "  +4.63%  "
$ws.Range("E34").Value = "  -1.80%  "
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("E36").Value = "  -0.78%  "
$ws.Range("D37").Value = "'0.581"
$ws.Range("E37").Value = "  -5.14%  "
$ws.Range("E38").Value = "  +0.37%  "
$ws.Range("E39").Value = "  +1.95%  "
$ws.Range("D40").Value = "'5.76"
$ws.Range("E40").Value = "  +3.66%  "
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("D42").Value = "'0.940"
$ws.Range("E42").Value = "  -15.71%  "
$ws.Range("E43").Value = "  +0.59%  "
$ws.Range("E44").Value = "  +0.57%  "
$ws.Range("D45").Value = "1.723.29"
$ws.Range("D46").Value = "'61.28"
$ws.Range("E46").Value = "  -1.71%  "
$ws.Range("D47").Value = "'85.87"
$ws.Range("E47").Value = "  -3.23%  "
$ws.Range("D48").Value = "0.0₆0101"
$ws.Range("E48").Value = "  -2.15%  "
$ws.Range("D49").Value = "'1.48"
$ws.Range("E49").Value = "  -2.52%  "
$ws.Range("D50").Value = "'0.0978"
$ws.Range("E50").Value = "  -2.67%  "
$ws.Range("E51").Value = "  -0.80%  "
